$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 287 (shifts COK and subsequent rows down by one)
$ws.Rows.Item(287).Insert()

# Populate the newly inserted row 287 with the Tahiti, French Polynesia entry
$ws.Cells.Item(287, 1).Value = "PPT"
$ws.Cells.Item(287, 2).Value = "Tahiti, French Polynesia"
$ws.Cells.Item(287, 3).Value = -17.5536994934
$ws.Cells.Item(287, 4).Value = -149.606994629
$ws.Cells.Item(287, 5).Value = "PF"
$ws.Cells.Item(287, 6).Value = "Oceania"
$ws.Cells.Item(287, 7).Value = "Tahiti"

# Match the style of column A used for colo codes in other rows (bold, bordered, centered)
$ws.Cells.Item(288, 1).Copy()
$ws.Cells.Item(287, 1).PasteSpecial(-4122)
